# Commit: "removed reference to 'predictive analytics'"
#   We have no code using machine-learning algorithms for predictive
#   analytics, so drop that bullet from the "TEAM 14's value proposition"
#   slide (slide 3), collapsing the now-redundant blank line with it so
#   "Competitive Intelligence" becomes the second bullet.

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(3)
$sh = $s.Shapes.Item(3)
$tr = $sh.TextFrame.TextRange

$full = $tr.Text
$idx0 = $full.IndexOf("Predictive Analytics")
$target = "Predictive Analytics" + [char]13 + [char]13
$len = $target.Length

$victim = $tr.Characters($idx0 + 1, $len)
$victim.Delete()
